$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new version-history row values first
$ws.Range("B8").Value = "_C3D-TEMPLATE_2025_FRA (Architecture v0001d)"
$ws.Range("C8").Value = "Ecrase le style dalle indice C pour le remplacer par le style dalle indice B (perte des jeux automatiques de propriétés)"

# Copy the formatting from the row above (row 6) so the new row matches
# the existing table styling (borders + left-aligned indent)
$ws.Range("B6:C6").Copy()
$ws.Range("B8:C8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply the values (PasteSpecial(xlPasteFormats) only carries formatting,
# but make sure the text content is exactly as intended)
$ws.Range("B8").Value = "_C3D-TEMPLATE_2025_FRA (Architecture v0001d)"
$ws.Range("C8").Value = "Ecrase le style dalle indice C pour le remplacer par le style dalle indice B (perte des jeux automatiques de propriétés)"

# Update the active-cell selection to match the new end-of-table location
$ws.Range("C12").Select()
